$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 1316
$ws.Range("F6").Value = 386
$ws.Range("F8").Value = 927
$ws.Range("F9").Value = 736
$ws.Range("F10").Value = 208
$ws.Range("F11").Value = 542
$ws.Range("F12").Value = 155
$ws.Range("F15").Value = 3045
$ws.Range("F16").Value = 2673
$ws.Range("F19").Value = 15
$ws.Range("F21").Value = 250
$ws.Range("F22").Value = 21
$ws.Range("F23").Value = 5452
$ws.Range("F24").Value = 598
$ws.Range("F26").Value = 37
$ws.Range("F27").Value = 67
$ws.Range("F28").Value = 380
$ws.Range("F31").Value = 84
$ws.Range("F32").Value = 305

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 20
$ws.Range("F7").Value = 21
$ws.Range("F9").Value = 42
$ws.Range("F16").Value = 993
$ws.Range("F23").Value = 332

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1789
$ws.Range("F5").Value = 2523
$ws.Range("F6").Value = 1083
$ws.Range("F9").Value = 1385
$ws.Range("F10").Value = 385

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1789
$ws.Range("F5").Value = 2523
$ws.Range("F9").Value = 1083
$ws.Range("F10").Value = 1385
$ws.Range("F11").Value = 385
$ws.Range("F14").Value = 1316
$ws.Range("F15").Value = 386
$ws.Range("F16").Value = 927
$ws.Range("F17").Value = 736
$ws.Range("F19").Value = 208
$ws.Range("F20").Value = 542
$ws.Range("F21").Value = 155
$ws.Range("F22").Value = 20
$ws.Range("F23").Value = 3045
$ws.Range("F24").Value = 2673
$ws.Range("F26").Value = 15
$ws.Range("F28").Value = 42
$ws.Range("F29").Value = 250
$ws.Range("F30").Value = 21
$ws.Range("F31").Value = 5452
$ws.Range("F32").Value = 598
$ws.Range("F35").Value = 37
$ws.Range("F36").Value = 67
$ws.Range("F37").Value = 380
$ws.Range("F42").Value = 332
$ws.Range("F49").Value = 84
$ws.Range("F50").Value = 305
